$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.573.58"
$ws.Range("E2").Value = "  -1.67%  "
$ws.Range("D3").Value = "2.020.65"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -10.11%  "
$ws.Range("E6").Value = "  -2.71%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.22"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.370"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.29%  "
$ws.Range("E11").Value = "  -2.46%  "
$ws.Range("E12").Value = "  -1.79%  "
$ws.Range("D13").Value = "2.318.40"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.762"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.34%  "
$ws.Range("D18").Value = "2.016.44"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").Value = "36.725.96"
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.12%  "
$ws.Range("D21").Value = "0.0₃0798"
$ws.Range("E21").Value = "  -4.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "221.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("E25").Value = "  +1.73%  "
$ws.Range("E26").Value = "  -7.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.84%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.134"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.82%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.00%  "
$ws.Range("E31").Value = "  -2.92%  "
$ws.Range("E32").Value = "  -1.92%  "
$ws.Range("E33").Value = "  -4.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0600"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.66%  "
$ws.Range("E36").Value = "  -3.92%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.77"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.46%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.31"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.13%  "
$ws.Range("E40").Value = "  +3.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0968"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.02%  "
$ws.Range("D43").Value = "1.459.02"
$ws.Range("E43").Value = "  +1.52%  "
$ws.Range("E44").Value = "  -2.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +38.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.39%  "
$ws.Range("E47").Value = "  -6.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.43%  "

Write-Host "Applied cryptos update"